$wb = $excel.ActiveWorkbook

# The workbook's three sheets ("0142", "0495a", "0495b") each carry a
# vertical H:I key/value metadata table. A new "OBSERV" measure-observation
# key needs to be inserted right before the existing "stratification" key,
# pushing "stratification" (and its value) down one row while leaving every
# other row untouched.
#
# We do this with a genuine row insert/delete pair (rather than rewriting
# cell .Value directly) so that the "stratification" row's cells - including
# ones holding an explicit empty-string value - are relocated intact instead
# of being reconstructed (and potentially collapsed to a truly blank cell).
#   1. Insert a blank row at row 11  -> old row 11 becomes row 12,
#                                        everything else shifts down too.
#   2. Delete the new row 13 (the shifted-down duplicate of the old,
#                              already-blank, row 12) to undo that extra
#                              shift for every row below.
#   3. Put the new "OBSERV" label in H11 (now vacated).
foreach ($ws in $wb.Worksheets) {
    $ws.Rows("11:11").Insert()
    $ws.Rows("13:13").Delete()
    $ws.Range("H11").Value = "OBSERV"
}

# Match the saved selection/active cell on every sheet.
foreach ($ws in $wb.Worksheets) {
    [void]$ws.Activate()
    [void]$ws.Range("H12").Select()
}

# Restore sheet "0142" (the first sheet) as the active tab.
[void]$wb.Worksheets.Item(1).Activate()
